$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errorMsg = "no such element: Unable to locate element: {`"method`":`"id`",`"selector`":`"menuTrigger`"}
  (Session info: chrome=56.0.2924.87)
  (Driver info: chromedriver=2.25.426923 (0390b88869384d6eb0d5d09729679f934aab9eed),platform=Windows NT 10.0.14393 x86_64) (WARNING: The server did not provide any stacktrace information)
Command duration or timeout: 10.30 seconds
For documentation on this error, please visit: http://seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.0.1', revision: '1969d75', time: '2016-10-18 09:49:13 -0700'
System info: host: 'MQCSERVER', ip: '172.16.0.7', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '1.8.0_121'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities [{applicationCacheEnabled=false, rotatable=false, mobileEmulationEnabled=false, networkConnectionEnabled=false, chrome={chromedriverVersion=2.25.426923 (0390b88869384d6eb0d5d09729679f934aab9eed), userDataDir=C:\Users\admin\AppData\Local\Temp\scoped_dir4684_19381}, takesHeapSnapshot=true, pageLoadStrategy=normal, databaseEnabled=false, handlesAlerts=true, hasTouchScreen=false, version=56.0.2924.87, platform=XP, browserConnectionEnabled=false, nativeEvents=true, acceptSslCerts=true, locationContextEnabled=true, webStorageEnabled=true, browserName=chrome, takesScreenshot=true, javascriptEnabled=true, cssSelectorsEnabled=true}]
Session ID: c0b590a1ff96b6bcb3d2cb09d870ccc7
*** Element info: {Using=id, value=menuTrigger}"

# New row 4: a failing TEAM Workspace test result, inserted before the
# previously-existing pass rows (which shift down to rows 5-7).
$ws.Range("A4").Value = "TC001"
$ws.Range("B4").Value = "Creating the TEAM Workspace"
$ws.Range("C4").Value = "TEAM Workspace Should be created successfully"
$ws.Range("D4").Value = $errorMsg
$ws.Range("E4").Value = "Fail"
# The long multi-line error text auto-expands the row height; restore the
# default (non-custom) row height to match the original report layout.
$ws.Rows("4").AutoFit()

# Rows 5-7: duplicate the original "pass" result row two more times.
$ws.Range("A5").Value = "TC001"
$ws.Range("B5").Value = "Creating the TEAM Workspace"
$ws.Range("C5").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D5").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E5").Value = "Pass"

$ws.Range("A6").Value = "TC001"
$ws.Range("B6").Value = "Creating the TEAM Workspace"
$ws.Range("C6").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D6").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E6").Value = "Pass"

$ws.Range("A7").Value = "TC001"
$ws.Range("B7").Value = "Creating the TEAM Workspace"
$ws.Range("C7").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D7").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E7").Value = "Pass"
